# Update the Iraq MSME "Employment (% of total)" row (row 12) with refreshed
# figures. The source sheet stores these percentages as text (e.g. "78.4"),
# so a leading apostrophe is used to keep Excel from re-typing the entry as
# a number - this preserves the original text storage while only changing
# the displayed/stored figures, matching the source data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "'78.35"
$ws.Range("C12").Value = "'20.24"
$ws.Range("D12").Value = "'98.59"
